$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 219, shifting rows 219:328 down to 220:329
$ws.Rows.Item(219).Insert()

# Populate the new row 219 with the fresh record
$ws.Cells.Item(219, 1).Value = 4
$ws.Cells.Item(219, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(219, 3).Value = "Los Lagos"
$ws.Cells.Item(219, 4).Value = 44917
$ws.Cells.Item(219, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(219, 5).Value = 10
$ws.Cells.Item(219, 6).Value = 100112032
$ws.Cells.Item(219, 7).Value = "Zapallo italiano"
$ws.Cells.Item(219, 8).Value = "Sin especificar"
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 70
$ws.Cells.Item(219, 11).Value = 10000
$ws.Cells.Item(219, 12).Value = 11000
$ws.Cells.Item(219, 13).Value = 10500
$ws.Cells.Item(219, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(219, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(219, 16).Value = 210
$ws.Cells.Item(219, 17).Value = 50
$ws.Cells.Item(219, 18).Value = "Hortaliza"
